$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9:AM9").Copy()
$ws.Range("A24:AM24").PasteSpecial()

$ws.Range("A10:AM10").Copy()
$ws.Range("A25:AM25").PasteSpecial()

$ws.Range("A12:AM12").Copy()
$ws.Range("A26:AM26").PasteSpecial()

$ws.Range("A13:AM13").Copy()
$ws.Range("A27:AM27").PasteSpecial()

$ws.Range("A14:AM14").Copy()
$ws.Range("A28:AM28").PasteSpecial()

$ws.Range("A24").Value = "2017_03_04_A_6-7(0)"
$ws.Range("A25").Value = "2017_03_04_A_6-7(0)"
$ws.Range("A26").Value = "2017_05_08_A_4-5(0)"
$ws.Range("A27").Value = "2017_05_08_A_5-4(0)"
$ws.Range("A28").Value = "2017_05_08_A_5-4(0)"

$ws.Range("A15").Value = "2017_07_06_C_3-4"
$ws.Range("A16").Value = "2017_07_06_C_3-4"
$ws.Range("A17").Value = "2017_07_06_C_4-3"

$ws.Range("E16").Select()
Write-Host "done"
